$d = $word.ActiveDocument

# 1. "left-bottom" -> "left-middle" inside the <margin> example text.
#    Only touch the "-bottom" portion so the untouched "left" run keeps its
#    own formatting and the "-" run absorbs the replacement text while the
#    (now redundant) "bottom" run disappears.
$d.Content.Find.Execute("-bottom", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "-middle", 2)

# 2. Add an explicit footer distance (720 twips = 36 pt) to the section's
#    page margins.
$d.PageSetup.FooterDistance = 36

Write-Output "done"
